# Auto-generated: apply crypto price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    # Force Excel to store the string as TEXT (never as a Number/Date),
    # matching the inlineStr/"t" cell type used throughout this sheet,
    # then strip the quote-prefix formatting the leading apostrophe adds
    # so the cell keeps its original (unstyled) look.
    $Cell.Value = "'" + $Text
    $Cell.ClearFormats()
}

$ws.Range("D2").Value = "27.134.02"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.645.64"
$ws.Range("E3").Value = "  +0.22%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("D5") "218.52"
$ws.Range("E5").Value = "  -0.29%  "
Set-TextValue $ws.Range("D6") "0.509"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  +1.68%  "
Set-TextValue $ws.Range("D9") "0.0626"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "1.873.09"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.635.36"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.52%  "
Set-TextValue $ws.Range("D15") "0.539"
$ws.Range("E15").Value = "  +2.33%  "
Set-TextValue $ws.Range("D16") "67.53"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").Value = "27.101.13"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +0.81%  "
Set-TextValue $ws.Range("D19") "221.99"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("E20").Value = "  -0.52%  "
Set-TextValue $ws.Range("D21") "6.77"
$ws.Range("E21").Value = "  +3.66%  "
Set-TextValue $ws.Range("D22") "4.44"
$ws.Range("E22").Value = "  +1.48%  "
Set-TextValue $ws.Range("D23") "2.43"
$ws.Range("E23").Value = "  +1.60%  "
Set-TextValue $ws.Range("D24") "9.25"
$ws.Range("E24").Value = "  +0.44%  "
Set-TextValue $ws.Range("D25") "147.44"
$ws.Range("E25").Value = "  -0.16%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.62%  "
Set-TextValue $ws.Range("D27") "7.40"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("E28").Value = "  +0.95%  "
Set-TextValue $ws.Range("D29") "15.83"
$ws.Range("E29").Value = "  +0.44%  "
Set-TextValue $ws.Range("D30") "0.0507"
$ws.Range("E30").Value = "  -0.43%  "
Set-TextValue $ws.Range("D31") "1.19"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("E32").Value = "  -0.38%  "
Set-TextValue $ws.Range("D33") "3.03"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").Value = "1.273.94"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +3.20%  "
Set-TextValue $ws.Range("D38") "0.544"
$ws.Range("E38").Value = "  +2.38%  "
Set-TextValue $ws.Range("D39") "0.845"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  -0.46%  "
Set-TextValue $ws.Range("D41") "0.811"
Set-TextValue $ws.Range("D42") "5.38"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "1.783.12"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") "2.15"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "63.07"
$ws.Range("E45").Value = "  +3.36%  "
Set-TextValue $ws.Range("D46") "92.70"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("E49").Value = "  -0.33%  "
Set-TextValue $ws.Range("D50") "7.73"
$ws.Range("E50").Value = "  +1.74%  "
Set-TextValue $ws.Range("D51") "0.0973"
$ws.Range("E51").Value = "  +0.31%  "
